$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.298.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.11%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.831.23'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.69%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.18%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.28'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -8.21%  '

# Row 6
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.18%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5181'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.48%  '

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -8.85%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06722'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.86%  '

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -8.50%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7626'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -7.04%  '

# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07677'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.72%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.836.01'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.45%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.33'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.73%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.009'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.91%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.25%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.03'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.97%  '

# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.19%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007861'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -4.21%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.330.34'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.10%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.077.40'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.42%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.526'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.32%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.398'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -7.42%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.881'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.54%  '

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.91%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.14'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.28%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.642'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.01%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.90'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.05%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '110.70'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.96%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.173'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.55%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.104'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.62%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08707'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.66%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04822'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.22%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.123'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.00%  '

# Row 35
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.12%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6789'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -9.58%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.088'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.71%  '

# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.95%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.199'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -9.24%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4891'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -8.53%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '111.12'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.13%  '

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8901'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.42%  '

# Row 43
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.117'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.47%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.651'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -7.27%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4174'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -9.63%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1251'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.79%  '

# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.055'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.22%  '

# Row 49
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05872'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.37%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.22'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.23%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.07'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.47%  '
